$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# --- Row 148: CyberASAP -------------------------------------------------
# Values are written in this precise order so that newly created shared
# strings come out in the same order as the reference workbook.
$ws.Cells.Item(148, 3).Value = "https://ktn-uk.org/programme/cyberasap/"
$ws.Cells.Item(148, 2).Value = "CyberASAP"
$ws.Cells.Item(148, 1).Value = "Pre-seed"
$ws.Cells.Item(148, 8).Value = "CyberASAP is funded by the UK Government Department for Digital, Culture Media & Sport (DCMS) and delivered through KTN and Innovate UK."

# --- New "Comment" column header / row2 sample --------------------------
$ws.Cells.Item(1, 8).Value = "Comment"
$ws.Cells.Item(2, 8).Value = "x"

# --- Row 149: Digital Catapult ------------------------------------------
$ws.Cells.Item(149, 2).Value = "Digital Catapult"
$ws.Cells.Item(149, 3).Value = "https://www.digicatapult.org.uk/"
$ws.Cells.Item(149, 7).Value = "info@digicatapult.org.uk"
$ws.Hyperlinks.Add($ws.Cells.Item(149, 7), "mailto:info@digicatapult.org.uk")
$ws.Cells.Item(149, 1).Value = "Hub"
$ws.Cells.Item(149, 6).Value = "101 Euston Road; London; NW1 2RA"
$ws.Cells.Item(149, 5).Value = "AI; VR; AR; IoT; 5G"
$ws.Cells.Item(148, 6).Value = "London; Endinbugh"
$ws.Cells.Item(149, 8).Value = "Accelerating early adoption of advanced digital technology. Digital Catapult specialises in Future Networks, AI and Immersive."

# --- Reused (already existing) strings -----------------------------------
$ws.Cells.Item(148, 5).Value = "cybersecurity"
$ws.Cells.Item(148, 7).Value = "contact on web page"

# --- Markdown-link formulas (same pattern as the rest of column D) -------
$ws.Cells.Item(148, 4).Formula = '=CONCATENATE("* [",B148,"](",C148,")")'
$ws.Cells.Item(148, 4).Style = "Normal"
$ws.Cells.Item(149, 4).Formula = '=CONCATENATE("* [",B149,"](",C149,")")'
$ws.Cells.Item(149, 4).Style = "Normal"

# G148 keeps the default (unstyled) look, unlike G149 which became a hyperlink
$ws.Cells.Item(148, 7).Style = "Normal"

# --- Final selection, matching the saved workbook state -------------------
$ws.Range("B152").Select()
